$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-empty price columns on existing rows 54-58 ---
# I54/I55/I56 also need their style switched from the "empty" (s=7) look to
# the "has a P.ea value" look (s=6), matching I57/I58 which already had one.
# Easiest reliable way to move a cell to another cell's exact style via COM
# is Copy + PasteSpecial(xlPasteFormats).

$ws.Cells.Item(57, 9).Copy()
$ws.Cells.Item(54, 9).PasteSpecial(-4122)
$ws.Cells.Item(54, 9).Value = 0.5421
$ws.Cells.Item(54, 10).Formula = "=H54*I54"

$ws.Cells.Item(57, 9).Copy()
$ws.Cells.Item(55, 9).PasteSpecial(-4122)
$ws.Cells.Item(55, 9).Value = 0.36076
$ws.Cells.Item(55, 10).Formula = "=H55*I55"

$ws.Cells.Item(57, 9).Copy()
$ws.Cells.Item(56, 9).PasteSpecial(-4122)
$ws.Cells.Item(56, 9).Value = 0.24812
$ws.Cells.Item(56, 10).Formula = "=H56*I56"

$ws.Cells.Item(57, 10).Formula = "=H57*I57"

$ws.Cells.Item(58, 1).Value = 1
$ws.Cells.Item(58, 7).Formula = "=E58*F58"
$ws.Cells.Item(58, 10).Formula = "=H58*I58"

# --- Insert two new rows for the fuse parts (F4266-ND / BK-6013-ND) ---
# Inserting at row 59 pushes the old "Sum" row (and everything below) down by
# two, and Excel auto-copies the formatting from the row above (58).

$ws.Rows.Item(59).Insert()
$ws.Rows.Item(59).Insert()

$ws.Rows.Item(59).RowHeight = 12.1
$ws.Rows.Item(60).RowHeight = 12.1

# Row 59: F4266-ND fuse
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "F4266-ND"
$ws.Cells.Item(59, 3).Value = "FUSE BLADE 20A/32V MINI FAST-AC"
$ws.Cells.Item(59, 4).Value = "F4-1"
$ws.Cells.Item(59, 5).Value = 1
$ws.Cells.Item(59, 6).Value = 0.4
$ws.Cells.Item(59, 7).Formula = "=E59*F59"
$ws.Cells.Item(59, 9).Value = 0.204
$ws.Cells.Item(59, 10).Formula = "=H59*I59"

# Row 60: BK-6013-ND fuse clip
$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(60, 2).Value = "BK-6013-ND"
$ws.Cells.Item(60, 3).Value = "FUSE CLIP AUTO 0.110X0.032""BLADE"
$ws.Cells.Item(60, 4).Value = "F4"
$ws.Cells.Item(60, 5).Value = 1
$ws.Cells.Item(60, 6).Value = 0.47
$ws.Cells.Item(60, 7).Formula = "=E60*F60"
$ws.Cells.Item(60, 9).Value = 0.34
$ws.Cells.Item(60, 10).Formula = "=H60*I60"

# --- Fix the "Sum" row formula (now row 61) to include the two new rows ---
$ws.Cells.Item(61, 7).Formula = "=SUM(G3:G60)"

# --- Update selection / view to match the authored state ---
$ws.Range("A54").Select()
$excel.ActiveWindow.ScrollRow = 22
